# Shrink fit calc: add a second sheet ("Sheet2") with a full shrink-fit
# analysis (ring/cap/liner dimensions, aluminum interface pressure,
# liner pressure at cryo temp, fluid pressure / friction on the cap),
# and wire a couple of new helper cells on Sheet1 (DELTA Do, Delta Max,
# elastic modulus E) that Sheet2 consumes via Sheet1!P9.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Sheet1: new inputs / outputs used by the shrink-fit sheet
# ---------------------------------------------------------------------

# New labeled columns on the header row (K8 "DELTA Do", P8 "Delta Max")
$ws1.Range("K8").Value = "DELTA Do"
$ws1.Range("P8").Value = "Delta Max"

# New "E" (elastic modulus) input next to delta T / CTE block
$ws1.Range("G5").Value = "E"
$ws1.Range("H5").Value = 10000000

# Row 9 data changes: ring OD/ID numbers updated, F9 becomes a plain
# input (wall thickness) instead of a derived formula, E9's formula is
# rewritten in terms of F9, and new columns K9/P9 are added.
$ws1.Range("D9").Value = 2.875
$ws1.Range("F9").Value = 0.125
$ws1.Range("E9").Formula = "=D9+F9"
$ws1.Range("K9").Formula = "=E9-H9"
$ws1.Range("L9").Value = 0.125
$ws1.Range("P9").Formula = "=K9/2"

# ---------------------------------------------------------------------
# Sheet2: shrink-fit analysis
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Properties table
$ws2.Range("C2").Value = "Properties"
$ws2.Range("C5").Value = "Poisson's Ratio"
$ws2.Range("D2").Value = "Aluminum"
$ws2.Range("E2").Value = "PTFE"
$ws2.Range("C3").Value = "CTE (in/in/°F)"
$ws2.Range("C4").Value = "Elastic Modulus (psi)"
$ws2.Range("C6").Value = "Temperature Change"

$ws2.Range("D3").Formula = "=0.0000131"
$ws2.Range("E3").Formula = "=0.000086*(5/9)"
$ws2.Range("D4").Value = 10000000
$ws2.Range("E4").Value = 500000
$ws2.Range("D5").Value = 0.33
$ws2.Range("E5").Value = 0.46
$ws2.Range("D6").Formula = "=(9/5)*(293-77)"
$ws2.Range("E6").Formula = "=D6"

# Ring / Cap / Liner dimension tables
$ws2.Range("C9").Value = "Inner Radius"
$ws2.Range("C8").Value = "Ring Dimensions"
$ws2.Range("C10").Value = "Outer Radius"
$ws2.Range("F8").Value = "Cap Dimensions"
$ws2.Range("D8").Value = "(inches)"
$ws2.Range("G8").Value = "(inches)"
$ws2.Range("J8").Value = "(inches)"
$ws2.Range("I8").Value = "Liner Dimensions"

$ws2.Range("F9").Value = "Inner Radius"
$ws2.Range("F10").Value = "Outer Radius"
$ws2.Range("I9").Value = "Inner Radius"
$ws2.Range("I10").Value = "Outer Radius"

$ws2.Range("D9").Value = 1.5
$ws2.Range("D10").Value = 2
$ws2.Range("G10").Formula = "=D9"
$ws2.Range("G9").Formula = "=G10-0.125"
$ws2.Range("J9").Formula = "=D9"
$ws2.Range("J10").Formula = "=J9+0.125"

# Aluminum interface pressure / friction table
$ws2.Range("C17").Value = "Aluminum Interface Pressure (psi)"
$ws2.Range("C20").Value = "Liner Pressure at Cryo Temp (psi)"

$ws2.Range("C11").Value = "Interface Depth"
$ws2.Range("C12").Value = "Interface Area (in^2)"
$ws2.Range("D23").Value = "Total Force (lb)"
$ws2.Range("I11").Value = "Thickness"
$ws2.Range("I12").Value = "Cryo Thickness"
$ws2.Range("I13").Value = "Change in Thickness"

$ws2.Range("D11").Value = 0.125
$ws2.Range("D12").Formula = "=2*3.141592*D9*D11"
$ws2.Range("J11").Formula = "=0.125"
$ws2.Range("J12").Formula = "=J11-J11*E3*E6"
$ws2.Range("J13").Formula = "=J11-J12"

$ws2.Range("C23").Value = "Fluid Pressure on Cap (psi)"
$ws2.Range("D17").Value = "Normal Force (lb)"
$ws2.Range("F17").Value = "Friction Force (lb)"
$ws2.Range("E17").Value = "Friction Coefficient"

$ws2.Range("C18").Formula = "=(D4*Sheet1!P9)/(2*D9^3)*((D10^2-D9^2)*(D9^2-G9^2)/(D10^2-G9^2))"
$ws2.Range("D18").Formula = "=C18*D12"
$ws2.Range("E18").Value = 1.2
$ws2.Range("F18").Formula = "=E18*D18"

$ws2.Range("C21").Formula = "=(((Sheet1!P9-J13)/2)/(J12))*E4"

$ws2.Range("C24").Value = 45
$ws2.Range("D24").Formula = "=C24*3.141592*D9^2"

# Number formatting to match Sheet1's convention for derived dimension
# deltas (4-decimal display format already used elsewhere in the book).
$ws1.Range("G9,I9,J9,K9,M9,N9,P9").NumberFormat = "0.0000"

# ---------------------------------------------------------------------
# View state: Sheet1 selection moves to P9, Sheet2 becomes the active
# (selected) tab with D14 selected.
# ---------------------------------------------------------------------
$ws1.Range("P9").Select()
$ws2.Range("D14").Select()
$ws2.Activate()
